$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "999/CCC"
$ws.Range("C2").Value = "I150156"
$ws.Range("D2").Value = "LATIFA BADRANE"
$ws.Range("N2").Value = 0

# Replace row 3 with blank-space placeholder values (like old row 4 layout),
# but keep amounts at 1000 and remove the old row 3 content
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1000

# Delete old row 4 entirely
$ws.Rows.Item(4).Delete()
